$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Professional summary paragraph - neutralize "all Black and
# Asian-American voters" -> "50M voters"
# ---------------------------------------------------------------------------
$summaryPara = $d.Paragraphs(4)
$summaryRange = $d.Range($summaryPara.Range.Start, $summaryPara.Range.End)
$summaryRange.Find.Execute("affecting all Black and Asian-American voters, developed geospatial ML", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "affecting 50M voters, developed geospatial ML", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "Partner - Siege Analytics" bullet - split the run so "50M" is
# its own bold / colored run (matching the "23%"/"64%" runs already there).
# ---------------------------------------------------------------------------
$bulletPara = $d.Paragraphs(10)
$bulletRange = $d.Range($bulletPara.Range.Start, $bulletPara.Range.End)
$bulletRange.Find.Execute("all Black and Asian-American") | Out-Null
$bulletRange.Text = "50M"

$bulletPara2 = $d.Paragraphs(10)
$boldRange = $d.Range($bulletPara2.Range.Start, $bulletPara2.Range.End)
$boldRange.Find.Execute("50M") | Out-Null
$boldRange.Font.Bold = 1
$boldRange.Font.Color = 5258796

# ---------------------------------------------------------------------------
# Change 4: "Geospatial Demographic Classification System" project impact
# line - neutralize language and add "nationwide".
# ---------------------------------------------------------------------------
$impactPara = $d.Paragraphs(43)
$impactRange = $d.Range($impactPara.Range.Start, $impactPara.Range.End)
$impactRange.Find.Execute("affecting all Black and Asian-American voters, improved electoral", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "affecting 50M voters nationwide, improved electoral", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: relocate the "Software Engineer - Mautinoa Technologies" block
# (Heading3 + 4 paragraphs) from right before "KEY PROJECTS" to right before
# "Senior Analyst - Myers Research" (i.e. right after the Siege Analytics
# bullets).
# ---------------------------------------------------------------------------
$anchor = $d.Paragraphs(15)   # "Senior Analyst - Myers Research ..."
$mautinoaTexts = @(
    "Software Engineer - Mautinoa Technologies (Austin, TX) | 2016 - 2018",
    "Software Development",
    "• Conceived, architected and engineered econometric simulation software for humanitarian crises intervention measurement",
    "• Liaised with data and engineering directors at multinational NGOs (UNICEF, IFRC)",
    "• Geospatial analysis on populations and boundaries for impact assessment"
)

for ($i = 0; $i -lt $mautinoaTexts.Count; $i++) {
    $anchor.Range.InsertParagraphBefore()
}

for ($i = 0; $i -lt $mautinoaTexts.Count; $i++) {
    $newPara = $d.Paragraphs(15 + $i)
    if ($i -gt 0) {
        $newPara.Style = "Normal"
    }
    $newPara.Range.Text = $mautinoaTexts[$i]
}

# Delete the original block, now shifted down by 5 paragraphs (was 30-34).
$origStart = $d.Paragraphs(35)
$origEnd = $d.Paragraphs(39)
$delRange = $d.Range($origStart.Range.Start, $origEnd.Range.End)
$delRange.Delete()

Write-Output "Done. Final paragraph count: $($d.Paragraphs.Count)"
